# Add an "average_greenspace" column (T) that averages the yearly values (D:R)
# for each county, and flip the sign convention of the existing percent_change
# column (S) from (D-R)/D to (R-D)/D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) New header cell T1 = "average_greenspace", formatted like S1
# ------------------------------------------------------------------
$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("T1").Value = "average_greenspace"

# ------------------------------------------------------------------
# 2) New data column T2:T67 = AVERAGE(D:R) per row, formatted to match
#    the other numeric columns (Calibri 11, 2 decimal places)
# ------------------------------------------------------------------
$ws.Cells.Item(1, 21).Copy()              # a cell carrying the plain "Calibri 11" look
$ws.Range("T2:T67").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("T2:T67").NumberFormat = "0.00"

$ws.Range("T2").Formula = "=AVERAGE(D2:R2)"
$ws.Range("T3:T66").Formula = "=AVERAGE(D3:R3)"
$ws.Range("T67").Formula = "=AVERAGE(D67:R67)"

# ------------------------------------------------------------------
# 3) Flip sign convention of column S: (D-R)/D -> (R-D)/D
# ------------------------------------------------------------------
$ws.Range("S2").Formula = "=(R2-D2)/D2"
$ws.Range("S3:S66").Formula = "=(R3-D3)/D3"
$ws.Range("S67").Formula = "=(R67-D67)/D67"

# ------------------------------------------------------------------
# 4) Column widths: widen S slightly, add a width for the new column T
# ------------------------------------------------------------------
$ws.Columns("S").ColumnWidth = 16.140625
$ws.Columns("T").ColumnWidth = 12.85546875

# ------------------------------------------------------------------
# 5) Selection ends up on the newly added T2 cell
# ------------------------------------------------------------------
$ws.Range("T2").Select()
